$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.626.71"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.848.97"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("D4").Value = "'1.030"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'321.68"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").Value = "'1.027"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.4375"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'0.3785"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("D9").Value = "'0.07378"
$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("D10").Value = "'0.8815"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").Value = "'21.50"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").Value = "1.860.23"
$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").Value = "'5.496"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").Value = "'6.713"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").Value = "'0.07132"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").Value = "'85.21"
$ws.Range("E16").Value = "  +3.00%  "

$ws.Range("D17").Value = "'1.032"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").Value = "'0.000009063"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").Value = "'1.026"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").Value = "'15.44"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").Value = "27.672.43"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "'5.285"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").Value = "'11.29"
$ws.Range("E23").Value = "  +0.93%  "

$ws.Range("D24").Value = "2.081.92"
$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D25").Value = "'2.028"
$ws.Range("E25").Value = "  +5.58%  "

$ws.Range("D26").Value = "'157.43"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").Value = "'18.68"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").Value = "'1.991"
$ws.Range("E28").Value = "  +2.71%  "

$ws.Range("D29").Value = "'5.328"
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("D30").Value = "'117.89"
$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("D31").Value = "'0.09008"
$ws.Range("E31").Value = "  -0.67%  "

$ws.Range("D32").Value = "'0.7717"
$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("D33").Value = "'1.208"
$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("D34").Value = "'2.992"
$ws.Range("E34").Value = "  +4.05%  "

$ws.Range("D35").Value = "'4.546"
$ws.Range("E35").Value = "  +1.29%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").Value = "'1.142"
$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("D38").Value = "'0.01970"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").Value = "'0.05271"
$ws.Range("E39").Value = "  +0.34%  "

$ws.Range("D40").Value = "'2.842"
$ws.Range("E40").Value = "  +2.19%  "

$ws.Range("D41").Value = "'0.5172"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").Value = "'0.1666"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "'6.841"
$ws.Range("E43").Value = "  +3.23%  "

$ws.Range("D44").Value = "'8.773"
$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("D45").Value = "'110.05"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("D46").Value = "'10.70"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").Value = "'0.06599"
$ws.Range("E47").Value = "  +4.31%  "

$ws.Range("D48").Value = "'1.029"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("D49").Value = "'1.700"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").Value = "'0.4691"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").Value = "'1.894"
$ws.Range("E51").Value = "  +0.14%  "
